# v2 Add more filters and to export Excel
#
# - Update the visit counters embedded in B1's free-text blob
#   (За сутки / За последние 30 дней / Всего).
# - Add a new column F1 with the "parnaya" (sauna type) filter text
#   that is already summarized at the start of C1 ("Вид парной: ...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the counters inside B1, leaving the rest of the text intact ---
$b1 = $ws.Range("B1").Value2
$b1 = $b1 -replace [regex]::Escape("За сутки: 226"), "За сутки: 262"
$b1 = $b1 -replace [regex]::Escape("За последние 30 дней: 336"), "За последние 30 дней: 372"
$b1 = $b1 -replace [regex]::Escape("Всего: 2466"), "Всего: 2502"
$ws.Range("B1").Value = $b1

# Re-setting a multi-line value makes the host auto-size the row; put the
# row back to its (default, non-custom) height so only the cell content
# changes, matching the source edit.
$ws.Rows.Item(1).EntireRow.AutoFit()

# --- Add the new F1 filter column ---
$ws.Range("F1").Value = "финская парная, русская баня, русская на дровах"
